$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value2 = "34.587.12"
$ws.Range("E2").Value2 = "  +1.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value2 = "1.799.59"
$ws.Range("E3").Value2 = "  +0.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value2 = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value2 = "227.56"
$ws.Range("E5").Value2 = "  +0.45%  "

# Row 6 - XRP
$ws.Range("D6").Value2 = "0.557"
$ws.Range("E6").Value2 = "  +1.71%  "

# Row 7 - USDC
$ws.Range("E7").Value2 = "  -0.08%  "

# Row 8 - Solana
$ws.Range("D8").Value2 = "32.73"
$ws.Range("E8").Value2 = "  +2.28%  "

# Row 9 - Cardano
$ws.Range("D9").Value2 = "0.298"
$ws.Range("E9").Value2 = "  +2.09%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value2 = "0.0697"
$ws.Range("E10").Value2 = "  +0.54%  "

# Row 11 - TRON
$ws.Range("E11").Value2 = "  +0.32%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value2 = "2.060.19"
$ws.Range("E12").Value2 = "  +0.92%  "

# Row 13 - Chainlink
$ws.Range("D13").Value2 = "11.18"
$ws.Range("E13").Value2 = "  +0.93%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value2 = "1.796.24"
$ws.Range("E14").Value2 = "  +0.86%  "

# Row 15 - Polygon
$ws.Range("E15").Value2 = "  +2.40%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value2 = "34.585.20"
$ws.Range("E16").Value2 = "  +1.23%  "

# Row 17 - Polkadot
$ws.Range("D17").Value2 = "4.32"
$ws.Range("E17").Value2 = "  +2.92%  "

# Row 18 - Litecoin
$ws.Range("D18").Value2 = "68.96"
$ws.Range("E18").Value2 = "  +1.43%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value2 = "0.0₃0805"
$ws.Range("E19").Value2 = "  +0.01%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value2 = "247.53"
$ws.Range("E20").Value2 = "  +0.13%  "

# Row 21 - Avalanche
$ws.Range("D21").Value2 = "11.30"
$ws.Range("E21").Value2 = "  +2.78%  "

# Row 22 - Dai
$ws.Range("E22").Value2 = "  -0.14%  "

# Row 23 - Uniswap
$ws.Range("E23").Value2 = "  +2.18%  "

# Row 24 - Monero
$ws.Range("D24").Value2 = "168.95"
$ws.Range("E24").Value2 = "  +3.81%  "

# Row 25 - Toncoin
$ws.Range("E25").Value2 = "  +1.51%  "

# Row 26 - Cosmos
$ws.Range("E26").Value2 = "  +1.48%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value2 = "16.62"
$ws.Range("E27").Value2 = "  +1.87%  "

# Row 28 - Stellar
$ws.Range("E28").Value2 = "  +2.27%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value2 = "  -0.04%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("E30").Value2 = "  +10.24%  "

# Row 31 & 32 swap - Hedera / Filecoin
$ws.Range("B31").Value2 = "Filecoin"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value2 = "3.83"
$ws.Range("E31").Value2 = "  +2.26%  "

$ws.Range("B32").Value2 = "Hedera"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value2 = "0.0527"
$ws.Range("E32").Value2 = "  +1.03%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value2 = "  +0.78%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value2 = "  +2.61%  "

# Row 35 - Maker
$ws.Range("D35").Value2 = "1.434.19"
$ws.Range("E35").Value2 = "  -0.70%  "

# Row 36 - RenderToken
$ws.Range("E36").Value2 = "  +8.09%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value2 = "0.675"
$ws.Range("E37").Value2 = "  +3.11%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value2 = "  +2.58%  "

# Row 39 - VeChain
$ws.Range("E39").Value2 = "  +0.33%  "

# Row 40 - Aave
$ws.Range("D40").Value2 = "85.18"
$ws.Range("E40").Value2 = "  +5.91%  "

# Row 41 & 42 swap - ARBITRUM / HuobiToken
$ws.Range("B41").Value2 = "HuobiToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value2 = "2.40"
$ws.Range("E41").Value2 = "  +1.35%  "

$ws.Range("B42").Value2 = "ARBITRUM"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value2 = "0.941"
$ws.Range("E42").Value2 = "  +1.71%  "

# Row 43 - MXToken
$ws.Range("D43").Value2 = "2.76"
$ws.Range("E43").Value2 = "  +3.33%  "

# Row 44 - InjectiveProtocol
$ws.Range("D44").Value2 = "13.88"
$ws.Range("E44").Value2 = "  +1.74%  "

# Row 45 - Kaspa
$ws.Range("E45").Value2 = "  +3.45%  "

# Row 46 - FraxShare
$ws.Range("E46").Value2 = "  +0.55%  "

# Row 47 - WEMIXToken
$ws.Range("E47").Value2 = "  +0.46%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value2 = "1.961.13"
$ws.Range("E48").Value2 = "  +0.92%  "

# Row 49 - Quant
$ws.Range("D49").Value2 = "106.14"
$ws.Range("E49").Value2 = "  +1.29%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value2 = "  -3.59%  "

# Row 51 - PaxDollar
$ws.Range("E51").Value2 = "  -0.10%  "
